$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header in B1 from "Gastos" to "Vendas"
$ws.Range("B1").Value = "Vendas"

# Add two new rows of data
$ws.Range("A8").Value = "açucar"
$ws.Range("B8").Value = 30

$ws.Range("A9").Value = "sal"
$ws.Range("B9").Value = 20

# Leave the selection on the next empty row, matching the author's final state
$ws.Range("B10").Select()
